{"js": "// Revert \"changed from wondows to windows\": find the misspelled-turned-\n// correct word \"windows\" (from \"Added this line in windows\") and change\n// it back to the original typo \"wondows\", collapsing the run/bookmark\n// split left behind by the earlier edit into a single run of text.\nconst body = context.document.body;\nconst results = body.search(\"windows\", { matchCase: false, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  const target = results.items[0];\n  target.insertText(\"wondows\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Revert \"changed from wondows to windows\": find the paragraph containing\n# the word \"windows\" (from \"Added this line in windows\") and restore the\n# original typo \"wondows\". The word is currently split across two runs\n# (\".../w\" + \"i\" + \"ndows\") with the \"_GoBack\" bookmark sitting between\n# them; a plain Find/Replace would either split a run around the bookmark\n# or silently drop the bookmark. Rebuilding the paragraph's text range via\n# InsertXML lets us land on the exact target shape: a single run holding\n# \"Added this line in wondows\" followed by the (still present) bookmark.\n$d = $word.ActiveDocument\n\n$search = $d.Content\n$search.Find.ClearFormatting()\n$search.Find.MatchCase = $true\n$search.Find.MatchWholeWord = $true\n$found = $search.Find.Execute(\"windows\")\n\nif ($found) {\n    # Expand the hit to its full enclosing paragraph (includes the trailing\n    # paragraph-mark character).\n    $paraRange = $search.Duplicate\n    $paraRange.Expand(4) | Out-Null  # wdParagraph\n\n    $fullText = $paraRange.Text\n    $markLen = 1\n    $bodyText = $fullText.Substring(0, $fullText.Length - $markLen)\n    $newText = $bodyText -replace \"windows\", \"wondows\"\n\n    # Range spanning only the paragraph's text (no trailing paragraph mark).\n    $bodyRange = $d.Range($paraRange.Start, $paraRange.End - $markLen)\n\n    function Esc([string]$s) {\n        return $s -replace '&', '&amp;' -replace '<', '&lt;' -replace '>', '&gt;'\n    }\n\n    if ($newText -ne $newText.Trim()) {\n        $spaceAttr = ' xml:space=\"preserve\"'\n    } else {\n        $spaceAttr = ''\n    }\n\n    $xml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n        '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        '<w:body><w:p>' +\n        '<w:r><w:t' + $spaceAttr + '>' + (Esc $newText) + '</w:t></w:r>' +\n        '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/>' +\n        '</w:p></w:body></w:document>' +\n        '</pkg:xmlData></pkg:part></pkg:package>'\n\n    $bodyRange.InsertXML($xml)\n}\n"}
